# Update filtered_output.xlsx: refresh the filtered stock-screener data on
# both the Neg_Change and Pos_Change sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Neg_Change ---------------------------------------------------
$wsNeg = $wb.Worksheets.Item("Neg_Change")

# Clear out the old data rows (2-7) before writing the new, shorter table.
$wsNeg.Range("A2:I7").Clear()

$negData = @(
    @("APOLLOHOSP", 7674.5, 7706.5, 7633.5, 7700,  210121,  423558,  -0.5039144579963075, "APOLLOHOSP"),
    @("KOTAKBANK",  2142,   2153.1, 2121.7, 2126,  2472939, 5152836, -0.5200819509877668, "KOTAKBANK"),
    @("MAXHEALTH",  1134,   1140,   1115,   1128,  3577437, 8494611, -0.5788580548302918, "MAXHEALTH"),
    @("AXISBANK",   1212,   1212,   1184,   1188,  8620156, 20128505, -0.5717438528097343, "AXISBANK")
)

$r = 2
foreach ($row in $negData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsNeg.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# --- Sheet 2: Pos_Change ----------------------------------------------------
$wsPos = $wb.Worksheets.Item("Pos_Change")

$wsPos.Range("A2:I6").Clear()

$posData = @(
    @("ETERNAL",    335.5,  341.6,  333.2,  337.85, 19836347, 13705095, 0.4473702663133674, "ETERNAL"),
    @("TCS",        2995,   3004.5, 2955.5, 2975,   3062943,  2112774,  0.4497258107114154, "TCS"),
    @("HINDUNILVR", 2540,   2559.4, 2508.3, 2530,   1291750,  818923,   0.5773766275950242, "HINDUNILVR"),
    @("NESTLEIND",  1181.7, 1192.2, 1175,   1176,   922441,   639965,   0.4413928886735993, "NESTLEIND"),
    @("JSWSTEEL",   1159.9, 1173,   1154.2, 1154.2, 1565661,  1013304,  0.5451049240899078, "JSWSTEEL")
)

$r = 2
foreach ($row in $posData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsPos.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
